{"js": "// Replace the date line and each division problem with its new value.\n// Every \"old\" string below is unique within the document, so a direct\n// search-and-replace (first match) is unambiguous.\nconst replacements = [\n  [\"2025-10-06 Monday\", \"2025-10-07 Tuesday\"],\n  [\"708\u00f76=\", \"349\u00f79=\"],\n  [\"555\u00f76=\", \"422\u00f79=\"],\n  [\"732\u00f75=\", \"715\u00f75=\"],\n  [\"257\u00f73=\", \"196\u00f79=\"],\n  [\"463\u00f72=\", \"649\u00f74=\"],\n  [\"558\u00f77=\", \"304\u00f72=\"],\n  [\"233\u00f73=\", \"580\u00f77=\"],\n  [\"705\u00f77=\", \"115\u00f77=\"],\n  [\"748\u00f72=\", \"777\u00f78=\"],\n  [\"694\u00f75=\", \"560\u00f75=\"],\n  [\"423\u00f75=\", \"572\u00f76=\"],\n  [\"983\u00f72=\", \"931\u00f78=\"],\n  [\"687\u00f78=\", \"784\u00f78=\"],\n  [\"962\u00f74=\", \"956\u00f72=\"],\n  [\"999\u00f74=\", \"191\u00f75=\"],\n  [\"205\u00f72=\", \"690\u00f78=\"],\n  [\"611\u00f72=\", \"864\u00f77=\"],\n  [\"254\u00f73=\", \"777\u00f72=\"],\n  [\"548\u00f73=\", \"390\u00f78=\"],\n  [\"758\u00f79=\", \"476\u00f73=\"],\n  [\"510\u00f76=\", \"659\u00f77=\"],\n  [\"772\u00f75=\", \"267\u00f78=\"],\n  [\"730\u00f72=\", \"577\u00f77=\"],\n  [\"724\u00f79=\", \"535\u00f73=\"],\n  [\"324\u00f72=\", \"962\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  // Old text is unique in the document, so only the first hit matters.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and each division problem with its new value.\n# Every \"old\" string is unique within the document, so Find/Replace across\n# the whole document body (wdReplaceAll) only ever touches one run each.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-10-06 Monday\", \"2025-10-07 Tuesday\"),\n    @(\"708\u00f76=\", \"349\u00f79=\"),\n    @(\"555\u00f76=\", \"422\u00f79=\"),\n    @(\"732\u00f75=\", \"715\u00f75=\"),\n    @(\"257\u00f73=\", \"196\u00f79=\"),\n    @(\"463\u00f72=\", \"649\u00f74=\"),\n    @(\"558\u00f77=\", \"304\u00f72=\"),\n    @(\"233\u00f73=\", \"580\u00f77=\"),\n    @(\"705\u00f77=\", \"115\u00f77=\"),\n    @(\"748\u00f72=\", \"777\u00f78=\"),\n    @(\"694\u00f75=\", \"560\u00f75=\"),\n    @(\"423\u00f75=\", \"572\u00f76=\"),\n    @(\"983\u00f72=\", \"931\u00f78=\"),\n    @(\"687\u00f78=\", \"784\u00f78=\"),\n    @(\"962\u00f74=\", \"956\u00f72=\"),\n    @(\"999\u00f74=\", \"191\u00f75=\"),\n    @(\"205\u00f72=\", \"690\u00f78=\"),\n    @(\"611\u00f72=\", \"864\u00f77=\"),\n    @(\"254\u00f73=\", \"777\u00f72=\"),\n    @(\"548\u00f73=\", \"390\u00f78=\"),\n    @(\"758\u00f79=\", \"476\u00f73=\"),\n    @(\"510\u00f76=\", \"659\u00f77=\"),\n    @(\"772\u00f75=\", \"267\u00f78=\"),\n    @(\"730\u00f72=\", \"577\u00f77=\"),\n    @(\"724\u00f79=\", \"535\u00f73=\"),\n    @(\"324\u00f72=\", \"962\u00f76=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
